# "Generate Report for Handoff"
#
# The localization status report records that the 4dfe909c… source file has
# now been handed off for translation (status flips from "In Translation" to
# "Ready for handoff", the priority becomes a machine-translation "mt" pass,
# and new handoff timestamps are recorded). Row 2 (the 4d0bdbdc… file) is
# untouched; only row 3 (the 4dfe909c… file) changes on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - summary columns for the 4dfe909c... row
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 18:16:34"

# ---------------------------------------------------------------------
# zh-cn sheet - Status / Priority / Latest Handoff Datetime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-05 18:16:29"

# ---------------------------------------------------------------------
# de-de sheet - Status / Priority / Latest Handoff Datetime
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-05 18:16:34"

# ---------------------------------------------------------------------
# The longer "Ready for handoff" status text no longer fits the old Status
# column width, so Excel widens the Status column on every sheet (E/F on
# Overview, C on the language sheets) to fit the new text.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
